# ---------------------------------------------------------------------
# Applies the "漫展信息" (convention info) update to the 展览 and 全部类型
# sheets:
#   - a handful of "想去人数" (interest count) bumps
#   - the 中秋嘉年华 event gets marked cancelled / ticket price -> "不可售"
#   - a brand-new 景德镇 event is inserted, the old 萌卡动漫展 row shifts
#     down one slot (its count also increases), and the JMG row shifts
#     down another slot (its count also increases)
# ---------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

function Set-TextCell($ws, $row, $col, $text) {
    # Forces the cell to be stored as text even when the value looks like
    # a date (e.g. "2024-10-02"), then restores the default ("Normal")
    # cell style so no stray direct formatting is left behind.
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

function Set-RowNumberStyle($ws, $row) {
    # Matches the bold / centered / bordered look used by every cell in
    # column A (style index 1 in the original workbook).
    $cell = $ws.Cells.Item($row, 1)
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop
    $cell.Borders.LineStyle = 1
}

function Apply-SheetChanges($ws, $rOffset) {
    # $rOffset is how many extra rows this sheet has above the common
    # "展览" rows before the shared tail section starts (0 for 展览,
    # 2 for 全部类型).
    $rAud    = 9  + $rOffset    # 南昌·Aud中秋动漫嘉年华 row
    $rMengka = 12 + $rOffset    # originally 南昌·萌卡动漫展 row
    $rJMG    = 13 + $rOffset    # originally 江西·JMG... row

    # ---- simple "想去人数" bumps (rows before 中秋嘉年华 row) ----
    if ($rOffset -eq 0) {
        $ws.Cells.Item(3, 6).Value = 1047
        $ws.Cells.Item(5, 6).Value = 2847
        $ws.Cells.Item(7, 6).Value = 241
    } else {
        $ws.Cells.Item(4, 6).Value = 1047
        $ws.Cells.Item(6, 6).Value = 2847
        $ws.Cells.Item(8, 6).Value = 241
    }

    # ---- 中秋嘉年华 marked as cancelled / not sellable ----
    $ws.Cells.Item($rAud, 3).Value = "南昌·Aud中秋动漫嘉年华（取消）"
    $ws.Cells.Item($rAud, 7).Value = "不可售"

    # ---- more "想去人数" bumps right after that ----
    $ws.Cells.Item($rAud + 1, 6).Value = 86
    $ws.Cells.Item($rAud + 2, 6).Value = 111

    # ---- make room for a new row by writing the JMG row's data one
    #      slot further down (kept as plain cell writes - no Rows.Insert
    #      shift - so existing numeric values elsewhere on the sheet
    #      don't get needlessly re-serialized/rounded) ----
    $rJMGNew = $rJMG + 1
    Set-RowNumberStyle $ws $rJMGNew
    $ws.Cells.Item($rJMGNew, 1).Value = $rJMG
    Set-TextCell $ws $rJMGNew 2 "2024-10-03"
    $ws.Cells.Item($rJMGNew, 3).Value = "江西·JMG（广电）第二届UP动漫游戏博览会"
    $ws.Cells.Item($rJMGNew, 4).Value = "怀玉山大道1315号 南昌绿地国际博览中心"
    $ws.Cells.Item($rJMGNew, 5).Value = "2024.10.03 09:00-10.05 18:00"
    $ws.Cells.Item($rJMGNew, 6).Value = 894
    $ws.Cells.Item($rJMGNew, 7).Value = 19.9
    $ws.Cells.Item($rJMGNew, 8).Value = "https://show.bilibili.com/platform/detail.html?id=90599"
    $ws.Cells.Item($rJMGNew, 9).Value = "//i2.hdslb.com/bfs/openplatform/202408/2LP6dm961723428231240.jpeg"

    # the row that used to hold JMG's data now gets the old 南昌·萌卡动漫展
    # info, with its "想去人数" bumped from 2667 to 2678. Its index number
    # keeps the value that used to belong to the JMG row (the source
    # data was generated this way - only the brand-new bottom row gets
    # a freshly incremented number).
    Set-RowNumberStyle $ws $rJMG
    $ws.Cells.Item($rJMG, 1).Value = $rJMG - 1
    Set-TextCell $ws $rJMG 2 "2024-10-02"
    $ws.Cells.Item($rJMG, 3).Value = "南昌·萌卡动漫展"
    $ws.Cells.Item($rJMG, 4).Value = "八一桥街道青山南路118号蓝海购物广场F1 蓝海展览馆"
    $ws.Cells.Item($rJMG, 5).Value = "2024.10.02 09:00-10.03 17:00"
    $ws.Cells.Item($rJMG, 6).Value = 2678
    $ws.Cells.Item($rJMG, 7).Value = 65
    $ws.Cells.Item($rJMG, 8).Value = "https://show.bilibili.com/platform/detail.html?id=89738"
    $ws.Cells.Item($rJMG, 9).Value = "//i0.hdslb.com/bfs/openplatform/202407/uqTvacSV1721621530709.jpeg"

    # row $rMengka (unchanged row number) becomes the brand-new
    # 景德镇 event
    $ws.Cells.Item($rMengka, 3).Value = "【大会员提前抢】景德镇·第十六届瓷都ACG内场—花玲&宴宁"
    $ws.Cells.Item($rMengka, 4).Value = "迎宾大道与寺山路交叉口东200米 陶博城"
    $ws.Cells.Item($rMengka, 5).Value = "2024.10.02 09:00-10.02 17:00"
    $ws.Cells.Item($rMengka, 6).Value = 3
    $ws.Cells.Item($rMengka, 7).Value = 188
    $ws.Cells.Item($rMengka, 8).Value = "https://show.bilibili.com/platform/detail.html?id=91080"
    $ws.Cells.Item($rMengka, 9).Value = "//i1.hdslb.com/bfs/openplatform/202408/z5XgI7ZE1724229859134.jpeg"
}

# ---------- Sheet: 展览 ----------
$ws1 = $wb.Worksheets.Item("展览")
Apply-SheetChanges $ws1 0

# ---------- Sheet: 全部类型 ----------
$ws4 = $wb.Worksheets.Item("全部类型")
Apply-SheetChanges $ws4 2
